$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B1 = 0, styled: bold font, thin box border, centered horizontally, top vertically
$ws.Range("B1").Value = 0
$ws.Range("B1").Font.Bold = $true
$ws.Range("B1").HorizontalAlignment = -4108   # xlCenter
$ws.Range("B1").VerticalAlignment = -4160     # xlTop
$ws.Range("B1").Borders.LineStyle = 1         # xlContinuous
$ws.Range("B1").Borders.Weight = 2            # xlThin

# A2 = 0, with the same style as B1 (copy formatting to avoid creating a
# redundant, unused style entry)
$ws.Range("A2").Value = 0
$ws.Range("B1").Copy()
$ws.Range("A2").PasteSpecial(-4122)           # xlPasteFormats

# B2 = "disconnected_elements" (default style, becomes a shared string)
$ws.Range("B2").Value = "disconnected_elements"

$excel.CutCopyMode = $false
